$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19.0299024646938
$ws.Range("C2").Value = 13.02988195902424
$ws.Range("D2").Value = 5.411078095326297
$ws.Range("E2").Value = 16.46589187120819
$ws.Range("F2").Value = 33.90383578391599
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("N2").Value = 17.44316607581648
$ws.Range("B3").Value = 18.2226977188296
$ws.Range("C3").Value = 12.24552864516589
$ws.Range("D3").Value = 5.420294845771187
$ws.Range("E3").Value = 15.52545175275263
$ws.Range("F3").Value = 33.22637460469882
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("N3").Value = 17.50383966367419
$ws.Range("B4").Value = 17.7152809864411
$ws.Range("C4").Value = 11.74112660364019
$ws.Range("D4").Value = 5.427424644875193
$ws.Range("E4").Value = 14.9245583793404
$ws.Range("F4").Value = 32.81680671120823
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("N4").Value = 17.5430766301334
$ws.Range("B5").Value = 17.50588947413828
$ws.Range("C5").Value = 11.53001738689059
$ws.Range("D5").Value = 5.430695470240081
$ws.Range("E5").Value = 14.67406152717285
$ws.Range("F5").Value = 32.65174673131767
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("N5").Value = 17.55956378142664
$ws.Range("B6").Value = 17.47097323410299
$ws.Range("C6").Value = 11.49463281516136
$ws.Range("D6").Value = 5.431260524904916
$ws.Range("E6").Value = 14.63213578358176
$ws.Range("F6").Value = 32.62445671002457
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("N6").Value = 17.56233149876924
$ws.Range("B7").Value = 17.71246715610447
$ws.Range("C7").Value = 11.73830177068373
$ws.Range("D7").Value = 5.427467282885502
$ws.Range("E7").Value = 14.92120248101157
$ws.Range("F7").Value = 32.81457288765623
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("N7").Value = 17.5432969673107
$ws.Range("B8").Value = 18.75422260181159
$ws.Range("C8").Value = 12.76426290818101
$ws.Range("D8").Value = 5.413948140114115
$ws.Range("E8").Value = 16.14664389643332
$ws.Range("F8").Value = 33.66906042346503
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("N8").Value = 17.4636741721765
$ws.Range("B9").Value = 20.69019709301865
$ws.Range("C9").Value = 14.58935103793376
$ws.Range("D9").Value = 5.399302330710603
$ws.Range("E9").Value = 18.43677942462618
$ws.Range("F9").Value = 35.3851708306037
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("N9").Value = 17.32331052248272
$ws.Range("B10").Value = 22.03211809840238
$ws.Range("C10").Value = 15.83207327609944
$ws.Range("D10").Value = 5.396039689456074
$ws.Range("E10").Value = 20.1034978793817
$ws.Range("F10").Value = 36.65709078807711
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("N10").Value = 17.22986056205369
$ws.Range("B11").Value = 22.62258163076735
$ws.Range("C11").Value = 16.39552342057327
$ws.Range("D11").Value = 5.396238391367107
$ws.Range("E11").Value = 20.82061698643642
$ws.Range("F11").Value = 37.23536968867223
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("N11").Value = 17.18946244811672
$ws.Range("B12").Value = 22.8431290553124
$ws.Range("C12").Value = 16.60418592823392
$ws.Range("D12").Value = 5.396559966279663
$ws.Range("E12").Value = 21.08630846517839
$ws.Range("F12").Value = 37.45408400282794
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("N12").Value = 17.17446992867203
$ws.Range("B13").Value = 22.79576830045271
$ws.Range("C13").Value = 16.55945556941721
$ws.Range("D13").Value = 5.39647968589023
$ws.Range("E13").Value = 21.02934728594566
$ws.Range("F13").Value = 37.40699525175848
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("N13").Value = 17.17768522870712
$ws.Range("B14").Value = 22.64078820047841
$ws.Range("C14").Value = 16.41278435378999
$ws.Range("D14").Value = 5.396259886268408
$ws.Range("E14").Value = 20.84259291084896
$ws.Range("F14").Value = 37.25337003105781
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("N14").Value = 17.18822287521207
$ws.Range("B15").Value = 22.54545669807582
$ws.Range("C15").Value = 16.32233208624714
$ws.Range("D15").Value = 5.396157458915321
$ws.Range("E15").Value = 20.7274377651631
$ws.Range("F15").Value = 37.15922896241352
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("N15").Value = 17.19471730391493
$ws.Range("B16").Value = 21.9931122596022
$ws.Range("C16").Value = 15.79459104202971
$ws.Range("D16").Value = 5.396060940770883
$ws.Range("E16").Value = 20.05580935891375
$ws.Range("F16").Value = 36.61927326808261
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("N16").Value = 17.23254331881922
$ws.Range("B17").Value = 21.64901578703598
$ws.Range("C17").Value = 15.46480566260421
$ws.Range("D17").Value = 5.396435876056854
$ws.Range("E17").Value = 19.63329498098246
$ws.Range("F17").Value = 36.28778456046521
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("N17").Value = 17.25629055856476
$ws.Range("B18").Value = 21.44922834138448
$ws.Range("C18").Value = 15.2836530714565
$ws.Range("D18").Value = 5.396809678362816
$ws.Range("E18").Value = 19.38640682776752
$ws.Range("F18").Value = 36.09710311631636
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("N18").Value = 17.27014810529746
$ws.Range("B19").Value = 21.38126794783762
$ws.Range("C19").Value = 15.22188702843998
$ws.Range("D19").Value = 5.39696326990292
$ws.Range("E19").Value = 19.30214897137628
$ws.Range("F19").Value = 36.0325454808209
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("N19").Value = 17.27487412402039
$ws.Range("B20").Value = 21.68584056545765
$ws.Range("C20").Value = 15.49812785968368
$ws.Range("D20").Value = 5.396379563967623
$ws.Range("E20").Value = 19.67867249184244
$ws.Range("F20").Value = 36.32307558604352
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("N20").Value = 17.25374204219344
$ws.Range("B21").Value = 22.68639360547148
$ws.Range("C21").Value = 16.45599272683964
$ws.Range("D21").Value = 5.396317727404059
$ws.Range("E21").Value = 20.89760607006963
$ws.Range("F21").Value = 37.29850244377744
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("N21").Value = 17.18511941245413
$ws.Range("B22").Value = 23.32248186787035
$ws.Range("C22").Value = 17.05461411943589
$ws.Range("D22").Value = 5.397715049526521
$ws.Range("E22").Value = 21.66007742077278
$ws.Range("F22").Value = 37.93434744083881
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("N22").Value = 17.14205127403829
$ws.Range("B23").Value = 22.98467222003882
$ws.Range("C23").Value = 16.73761883403753
$ws.Range("D23").Value = 5.396836316537512
$ws.Range("E23").Value = 21.25624525779039
$ws.Range("F23").Value = 37.59520630896871
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("N23").Value = 17.16487407481469
$ws.Range("B24").Value = 21.66919818016649
$ws.Range("C24").Value = 15.48307098709758
$ws.Range("D24").Value = 5.396404530088661
$ws.Range("E24").Value = 19.65816971602738
$ws.Range("F24").Value = 36.3071208209746
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("N24").Value = 17.25489358821267
$ws.Range("B25").Value = 20.17962238864113
$ws.Range("C25").Value = 14.1163163754722
$ws.Range("D25").Value = 5.401966546699272
$ws.Range("E25").Value = 17.78629376902386
$ws.Range("F25").Value = 34.91800615943511
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("N25").Value = 17.35958798893854
